# Trade #89 closed at 2026-02-16 21:38:07 - momentum DOWN +0.000%
#
# This script applies the following changes to the workbook:
#  1. "Summary" sheet: refresh OVERALL and leadlag aggregate stats.
#  2. "leadlag" sheet: close out trade #58 (row 47) that was previously OPEN.
#  3. "momentum" sheet: append newly opened trade #89 as a new row.
#  4. "All Trades" sheet: append the now-CLOSED leadlag trade #58 as a new row.
#  5. "Comparison" sheet: refresh leadlag aggregate stats.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 3).Value = 58
$summary.Cells.Item(2, 4).Value = "'67.2%"
$summary.Cells.Item(2, 5).Value = "'+15.4592%"
$summary.Cells.Item(2, 6).Value = "'+0.2665%"

$summary.Cells.Item(3, 4).Value = "'43.9%"
$summary.Cells.Item(3, 5).Value = "'+10.6084%"
$summary.Cells.Item(3, 6).Value = "'+0.1607%"

# ---------------------------------------------------------------------
# 2. leadlag sheet - close trade #58 in row 47
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Cells.Item(47, 7).Value = 68556.074761
$leadlag.Cells.Item(47, 8).Value = "CLOSED"
$leadlag.Cells.Item(47, 9).Value = 0.362
$leadlag.Cells.Item(47, 10).Value = 3.62
$leadlag.Cells.Item(47, 13).Value = "time_exit_5min"
$leadlag.Cells.Item(47, 14).Value = 5

# ---------------------------------------------------------------------
# 3. momentum sheet - append newly opened trade #89 as row 24
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Cells.Item(24, 1).Value = 89
$momentum.Cells.Item(24, 2).Value = "'2026-02-16"
$momentum.Cells.Item(24, 3).Value = "'21:38:07"
$momentum.Cells.Item(24, 4).Value = "momentum"
$momentum.Cells.Item(24, 5).Value = "DOWN"
$momentum.Cells.Item(24, 6).Value = 68443.59
# Column G (Exit Price) and M (Exit Reason) stay blank - trade is still OPEN.
$momentum.Cells.Item(24, 8).Value = "OPEN"
$momentum.Cells.Item(24, 9).Value = 0
$momentum.Cells.Item(24, 10).Value = 0
$momentum.Cells.Item(24, 11).Value = 0.9
$momentum.Cells.Item(24, 12).Value = "Downward momentum: -0.225% over 10 samples"
$momentum.Cells.Item(24, 14).Value = 0

# ---------------------------------------------------------------------
# 4. All Trades sheet - append the now-CLOSED leadlag trade #58 as row 59
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(59, 1).Value = 58
$allTrades.Cells.Item(59, 2).Value = "'2026-02-16"
$allTrades.Cells.Item(59, 3).Value = "'21:33:02"
$allTrades.Cells.Item(59, 4).Value = "leadlag"
$allTrades.Cells.Item(59, 5).Value = "DOWN"
$allTrades.Cells.Item(59, 6).Value = 68805.145
$allTrades.Cells.Item(59, 7).Value = 68556.074761
$allTrades.Cells.Item(59, 8).Value = "CLOSED"
$allTrades.Cells.Item(59, 9).Value = 0.362
$allTrades.Cells.Item(59, 10).Value = 3.62
$allTrades.Cells.Item(59, 11).Value = 0.75
$allTrades.Cells.Item(59, 12).Value = "Coinbase leading with -0.116% move"
$allTrades.Cells.Item(59, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(59, 14).Value = 5

# ---------------------------------------------------------------------
# 5. Comparison sheet - refresh leadlag aggregate stats
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Cells.Item(2, 3).Value = "'43.9%"
$comparison.Cells.Item(2, 4).Value = "'3.00"
$comparison.Cells.Item(2, 5).Value = "'+0.5486%"
$comparison.Cells.Item(2, 7).Value = "'1.76"
